# Update "paises" (countries) workbook: refresh COVID counters for a batch
# of countries, re-sort three country names into their correct alphabetical
# slot (Zambia, Belice, Santa Lucia), and bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "last updated" timestamp (A1) ---------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 15:39"

# --- Country-name shuffle (shared-string ordering fix) ---------------------
# Zambia moves up, ahead of Guinea/Tayikistan/Haiti/Finlandia (rows 91-95)
$ws.Range("A91").Value = "Zambia"
$ws.Range("A92").Value = "Guinea"
$ws.Range("A93").Value = "Tayikistan"
$ws.Range("A94").Value = "Haiti"
$ws.Range("A95").Value = "Finlandia"

# Belice moves up, ahead of Brunei/Barbados/Monaco/Seychelles (rows 186-190)
$ws.Range("A186").Value = "Belice"
$ws.Range("A187").Value = "Brunei"
$ws.Range("A188").Value = "Barbados"
$ws.Range("A189").Value = "Monaco"
$ws.Range("A190").Value = "Seychelles"

# Santa Lucia moves ahead of Timor Oriental (rows 202-203)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Data refresh (Casos totales, Nuevos casos, Casos activos, Recuperados,
#     Casos criticos, Muertes hoy, Muertes) for the affected countries -------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5098494
$ws.Range("C4").Value = 2970
$ws.Range("E4").Value = 2316325
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 164144

# Row 6 - India
$ws.Range("B6").Value = 2108705
$ws.Range("C6").Value = 21841
$ws.Range("D6").Value = 1443183
$ws.Range("E6").Value = 622724
$ws.Range("G6").Value = 220
$ws.Range("H6").Value = 42798

# Row 16 - Arabia Saudita
$ws.Range("B16").Value = 287262
$ws.Range("C16").Value = 1469
$ws.Range("D16").Value = 250440
$ws.Range("E16").Value = 33692
$ws.Range("G16").Value = 37
$ws.Range("H16").Value = 3130

# Row 41 - Kuwait
$ws.Range("B41").Value = 71199
$ws.Range("C41").Value = 472
$ws.Range("D41").Value = 62806
$ws.Range("E41").Value = 7919
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 474

# Row 45 - Paises Bajos
$ws.Range("B45").Value = 57987
$ws.Range("C45").Value = 486
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 6157

# Row 63 - Serbia
$ws.Range("B63").Value = 27863
$ws.Range("C63").Value = 255
$ws.Range("E63").Value = 13184
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 632

# Row 82 - Madagascar
$ws.Range("B82").Value = 12922
$ws.Range("C82").Value = 214
$ws.Range("D82").Value = 10604
$ws.Range("E82").Value = 2177
$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 141

# Row 91 - Zambia (newly positioned here, fresh data)
$ws.Range("B91").Value = 7903
$ws.Range("C91").Value = 417
$ws.Range("D91").Value = 6431
$ws.Range("E91").Value = 1269
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 203

# Row 92 - Guinea
$ws.Range("B92").Value = 7777
$ws.Range("D92").Value = 6800
$ws.Range("E92").Value = 927
$ws.Range("H92").Value = 50

# Row 93 - Tayikistan
$ws.Range("B93").Value = 7706
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 6484
$ws.Range("E93").Value = 1160
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 62

# Row 94 - Haiti
$ws.Range("B94").Value = 7599
$ws.Range("C94").Value = 17
$ws.Range("D94").Value = 4893
$ws.Range("E94").Value = 2529
$ws.Range("G94").Value = 6
$ws.Range("H94").Value = 177

# Row 95 - Finlandia
$ws.Range("B95").Value = 7568
$ws.Range("C95").Value = 14
$ws.Range("D95").Value = 6980
$ws.Range("E95").Value = 257
$ws.Range("H95").Value = 331

# Row 186 - Belice (newly positioned here, fresh data)
$ws.Range("B186").Value = 146
$ws.Range("C186").Value = 32
$ws.Range("D186").Value = 32
$ws.Range("E186").Value = 112
$ws.Range("H186").Value = 2

# Row 187 - Brunei
$ws.Range("B187").Value = 142
$ws.Range("D187").Value = 138
$ws.Range("E187").Value = 1
$ws.Range("H187").Value = 3

# Row 188 - Barbados
$ws.Range("B188").Value = 138
$ws.Range("D188").Value = 100
$ws.Range("E188").Value = 31
$ws.Range("H188").Value = 7

# Row 189 - Monaco
$ws.Range("B189").Value = 128
$ws.Range("D189").Value = 105
$ws.Range("E189").Value = 19
$ws.Range("H189").Value = 4

# Row 190 - Seychelles
$ws.Range("B190").Value = 126
$ws.Range("D190").Value = 125
$ws.Range("E190").Value = 1
$ws.Range("H190").Value = 0

# Row 212 - Montserrat
$ws.Range("D212").Value = 11
$ws.Range("E212").Value = 1
